$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF".
# Copy the formatting from H1 (bold/bordered/centered header style) onto the
# new header cells so they reuse the same cell style as the rest of row 1,
# then set their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Data rows 2-35: I = 1 (constant), J = same value as H (copied from column H)
for ($r = 2; $r -le 35; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
